$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 110
$ws.Range("AD2").Value = 18
$ws.Range("AK2").Value = 17.5
$ws.Range("AN2").Value = 10
$ws.Range("AO2").Value = 46
$ws.Range("F2").Value = 1.83
$ws.Range("J2").Value = 4.1
$ws.Range("P2").Value = 2.24
$ws.Range("S2").Value = 2.88
$ws.Range("T2").Value = 1.72
$ws.Range("W2").Value = 2.16
$ws.Range("X2").Value = 20
$ws.Range("AB4").Value = 8.800000000000001
$ws.Range("AC4").Value = 9.800000000000001
$ws.Range("AM4").Value = 140
$ws.Range("F4").Value = 1.6
$ws.Range("H4").Value = 5.6
$ws.Range("I4").Value = 6.6
$ws.Range("J4").Value = 3.95
$ws.Range("N4").Value = 3.9
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 1.39
$ws.Range("T4").Value = 1.84
$ws.Range("U4").Value = 1.98
$ws.Range("X4").Value = 20
$ws.Range("H5").Value = 1.06
$ws.Range("J5").Value = 1.03
$ws.Range("Q5").Value = 1.79
$ws.Range("R5").Value = 1.14
$ws.Range("S5").Value = 1.79
$ws.Range("F6").Value = 2.36
$ws.Range("G6").Value = 2.7
$ws.Range("H6").Value = 3.15
$ws.Range("I6").Value = 4.1
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3.5
$ws.Range("V6").Value = 1.32
$ws.Range("W6").Value = 1.53
$ws.Range("AC7").Value = 13
$ws.Range("AD7").Value = 11.5
$ws.Range("AE7").Value = 16
$ws.Range("AF7").Value = 980
$ws.Range("AH7").Value = 23
$ws.Range("AJ7").Value = 130
$ws.Range("AK7").Value = 55
$ws.Range("AM7").Value = 70
$ws.Range("F7").Value = 4.7
$ws.Range("G7").Value = 5.6
$ws.Range("H7").Value = 1.59
$ws.Range("I7").Value = 1.72
$ws.Range("J7").Value = 4.6
$ws.Range("K7").Value = 5.5
$ws.Range("N7").Value = 6.2
$ws.Range("O7").Value = 1.15
$ws.Range("P7").Value = 2.78
$ws.Range("Q7").Value = 1.45
$ws.Range("R7").Value = 1.71
$ws.Range("S7").Value = 2.12
$ws.Range("T7").Value = 1.55
$ws.Range("U7").Value = 2.42
$ws.Range("V7").Value = 2.38
$ws.Range("Y7").Value = 15
$ws.Range("Z7").Value = 14.5
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 9.800000000000001
$ws.Range("AH8").Value = 25
$ws.Range("AK8").Value = 140
$ws.Range("AL8").Value = 130
$ws.Range("G8").Value = 8.199999999999999
$ws.Range("K8").Value = 4.7
$ws.Range("L8").Value = 1.41
$ws.Range("N8").Value = 3.85
$ws.Range("P8").Value = 1.95
$ws.Range("Q8").Value = 1.94
$ws.Range("T8").Value = 2.04
$ws.Range("U8").Value = 1.84
$ws.Range("Y8").Value = 7.8
$ws.Range("F9").Value = 1.64
$ws.Range("G9").Value = 1.72
$ws.Range("I9").Value = 7.4
$ws.Range("P9").Value = 1.7
$ws.Range("Q9").Value = 2.24
$ws.Range("S9").Value = 4.3
$ws.Range("T9").Value = 2.16
$ws.Range("U9").Value = 1.74
$ws.Range("V9").Value = 1.15
$ws.Range("W9").Value = 2.38
$ws.Range("X9").Value = 13.5
$ws.Range("AE10").Value = 29
$ws.Range("AF10").Value = 19.5
$ws.Range("AK10").Value = 32
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 3.75
$ws.Range("O10").Value = 1.34
$ws.Range("T10").Value = 1.79
$ws.Range("U10").Value = 2.2
$ws.Range("O11").Value = 1.33
$ws.Range("AJ12").Value = 30
$ws.Range("AM12").Value = 370
$ws.Range("AO12").Value = 170
$ws.Range("M12").Value = 1.18
$ws.Range("I13").Value = 5.3
$ws.Range("J13").Value = 3.6
$ws.Range("T13").Value = 1.96
$ws.Range("U13").Value = 1.9
$ws.Range("V13").Value = 1.23
